$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes (D:G widened) ---
# ColumnWidth is expressed in "characters" and gets rounded by Excel to the
# nearest pixel (MDW=7) when stored, so we pick the input that lands on the
# pixel count closest to the target stored width.
$ws.Columns.Item(4).ColumnWidth = 46.57142857142857
$ws.Columns.Item(5).ColumnWidth = 38.142857142857146
$ws.Columns.Item(6).ColumnWidth = 50.714285714285715
$ws.Columns.Item(7).ColumnWidth = 41.857142857142854

# --- Row 32: Anteil erneuerbarer Energien am Brutto(-)Endenergieverbrauch ---
$ws.Range("F32").Value = "Anteil erneuerbarer Energien am Bruttoendenergieverbrauch"
$ws.Range("L32").Value = "Anteil erneuerbarer Energien am Bruttoendenergieverbrauch"

# --- Row 38: Verhältnis der Bruttoanlageinvestitionen zum BIP/Bruttoinlandsprodukt ---
$ws.Range("F38").Value = "Verhältnis der Bruttoanlageinvestitionen zum Bruttoinlandsprodukt"
$ws.Range("L38").Value = "Verhältnis der Bruttoanlageinvestitionen zum Bruttoinlandsprodukt"

# --- Row 39: Bruttoinlandsprodukt je Einwohner(in und Einwohner) ---
$ws.Range("F39").Value = "Bruttoinlandsprodukt je Einwohnerin und Einwohner"
$ws.Range("L39").Value = "Bruttoinlandsprodukt je Einwohnerin und Einwohner"

# --- Row 40: Erwerbstätigenquote insgesamt (20-64 Jahre) ---
$ws.Range("D40").Value = "a) Insgesamt (20-64-Jährige)"
$ws.Range("E40").Value = "a) total (20 to 64-year-olds)"
$ws.Range("F40").Value = "Erwerbstätigenquote (20- bis 64-Jährige)"
$ws.Range("G40").Value = "Employment rate (20 to 64-year-olds)"
$ws.Range("L40").Value = "Erwerbstätigenquote insgesamt (20- bis 64-Jährige)"

# --- Row 41: Erwerbstätigenquote Ältere (60-64 Jahre) ---
$ws.Range("D41").Value = "b) Ältere (60-64-Jährige)"
$ws.Range("F41").Value = "Erwerbstätigenquote (60- bis 64-Jährige)"
$ws.Range("G41").Value = "Employment rate (60 to 64-year-olds)"
$ws.Range("L41").Value = "Erwerbstätigenquote Ältere (60- bis 64-Jährige)"

# --- Row 44: Roll-out of broadband - share/Share of households ---
$ws.Range("G44").Value = "Roll-out of broadband – Share of households with access to gigabit broadband services"
$ws.Range("M44").Value = "Roll-out of broadband – Share of households with access to gigabit broadband services"

# --- Row 47: Anstieg der Siedlungs- und Verkehrsfläche (drop "in ha pro Tag") ---
$ws.Range("F47").Value = "Anstieg der Siedlungs- und Verkehrsfläche"
$ws.Range("G47").Value = "Expansion of settlement and transport area"
$ws.Range("L47").Value = "Anstieg der Siedlungs- und Verkehrsfläche"
$ws.Range("M47").Value = "Expansion of settlement and transport area"

# --- Row 53: Housing cost overload -> overburden ---
$ws.Range("G53").Value = "Housing cost overburden"
$ws.Range("M53").Value = "Housing cost overburden"

# --- Row 59: EMAS eco-management -> Eco-management EMAS ---
$ws.Range("E59").Value = "Eco-management EMAS"
$ws.Range("G59").Value = "Eco-management EMAS"
$ws.Range("M59").Value = "Eco-management EMAS"

# --- Row 64: Nitrogen input via the inflows into the Baltic Sea ---
$ws.Range("F64").Value = "Stickstoffeintrag über die Zuflüsse in die Ostsee"
$ws.Range("G64").Value = "Nitrogen input via the inflows into the Baltic Sea"
$ws.Range("M64").Value = "Nitrogen input via the inflows into the Baltic Sea"

# --- Row 65: Nitrogen input via the inflows into the North Sea ---
$ws.Range("F65").Value = "Stickstoffeintrag über die Zuflüsse in die Nordsee"
$ws.Range("G65").Value = "Nitrogen input via the inflows into the North Sea"
$ws.Range("M65").Value = "Nitrogen input via the inflows into the North Sea"

# --- Row 66: Share of sustainably fished (fish) stocks of fish ---
$ws.Range("G66").Value = "Share of sustainably fished stocks of fish in the North and Baltic Seas"
$ws.Range("M66").Value = "Share of sustainably fished stocks of fish in the North and Baltic Seas"

# --- Row 76: Studierende und Forschende aus Entwicklungs- und Schwellenländern / LDCs ---
$ws.Range("D76").Value = "Studierende und Forschende aus Entwicklungslenländern und LDCs"
$ws.Range("E76").Value = "Students and researchers from developing countries and LDCs"
$ws.Range("F76").Value = "Anzahl der Studierenden und Forschenden aus Entwicklungsländern sowie aus am wenigsten entwickelten Ländern pro Jahr"
$ws.Range("G76").Value = "Number of students and researchers from developing countries and least developed countries per year"
$ws.Range("L76").Value = "Anzahl der Studierenden und Forschenden aus Entwicklungsländern sowie aus am wenigsten entwickelten Ländern pro Jahr"
$ws.Range("M76").Value = "Number of students and researchers from developing countries and least developed countries per year"

Write-Host "Edits applied"
